# Add "PMID" column to the "studies" sheet (column H) and a "notes"
# column to the "counts" sheet (column F), per commit message:
# "added PMID and notes columns to data structure"

$wb = $excel.ActiveWorkbook

# --- studies sheet: add PMID header in column H ---
$wsStudies = $wb.Worksheets.Item("studies")
$wsStudies.Range("H1").Value = "PMID"

# --- counts sheet: add notes header in column F ---
$wsCounts = $wb.Worksheets.Item("counts")
$wsCounts.Range("F1").Value = "notes"

# --- update selection / active sheet state to match authored edit ---
# studies sheet is no longer the active tab; leave selection on the
# newly added column's first data cell
$wsStudies.Activate()
$wsStudies.Range("H2").Select()

# counts sheet becomes the active tab, selection on its new column's
# first data cell
$wsCounts.Activate()
$wsCounts.Range("F2").Select()
